$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 4273.066666666667, 4649, 3812, 0.09084692796071371),
    @(1, 4389.033333333334, 4724, 3965, 0.08979980945587158),
    @(2, 4558.133333333333, 4971, 3870, 0.09462287425994872),
    @(3, 4267.033333333334, 4628, 3367, 0.09210739135742188),
    @(4, 3778.8,             4051, 3490, 0.09511305491129557),
    @(5, 3998.833333333333, 4300, 3596, 0.09204138120015462),
    @(6, 4651.2,             5005, 4186, 0.09521447817484538),
    @(7, 4125.833333333333, 4511, 3745, 0.09506675402323404),
    @(8, 4290.333333333333, 4646, 3904, 0.09386487007141113),
    @(9, 4185.9,             4576, 3737, 0.09093227386474609)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
